# daily auto push: 2026-01-24 22:33 UTC
# Insert a new data row (2026/01/25, 日, 5, 18) right before the existing
# row 692 ("2026/12/29"), pushing every subsequent row down by one.
# The sheet's used range grows from A1:D733 to A1:D734.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 692:733 down to 693:734, opening up a blank row 692.
$ws.Rows.Item(692).Insert()

# Force column A to be treated as literal text so the date-looking
# string "2026/01/25" isn't auto-converted into a real Excel date
# serial (matches how every other date cell in this column is stored
# as plain text), then restore the default (unstyled) cell style so no
# stray number-format style lingers on the cell.
$ws.Range("A692").NumberFormat = "@"
$ws.Range("A692").Value = "2026/01/25"
$ws.Range("A692").Style = "Normal"

$ws.Range("B692").Value = "日"
$ws.Range("C692").Value = 5
$ws.Range("D692").Value = 18
